$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Reference cells used as style/value donors for type-changing cells ---
# "0" placeholder text (style 14, shared string "0")
$zeroRef = $ws.Range("C14")
# "***.*" placeholder text (style 14, shared string "***.*")
$naRef = $ws.Range("E14")
# numeric-style (165 "#,##0") donor cell for Count-type columns (style 15)
$numRef15 = $ws.Range("C16")
# numeric-style (167 "#,##0.0") donor cell for %Chg-type columns (style 16)
$numRef16 = $ws.Range("E16")

# --- Row 15 ---
$zeroRef.Copy($ws.Range("D15"))
$naRef.Copy($ws.Range("E15"))
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -81.081081081081

# --- Row 16 ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 148
$ws.Range("K16").Value = -2.702702702702
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -26.903553299492
$ws.Range("N16").Value = -79.831932773109

# --- Row 17 ---
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 194
$ws.Range("K17").Value = 5.670103092783
$ws.Range("L17").Value = -4.651162790697
$ws.Range("M17").Value = 31.410256410256
$ws.Range("N17").Value = -66.39344262295

# --- Row 18 ---
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 93
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = -38
$ws.Range("L18").Value = -31.111111111111
$ws.Range("M18").Value = 2.197802197802
$ws.Range("N18").Value = -87.13692946058

# --- Row 19 ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -13.157894736842
$ws.Range("I19").Value = 314
$ws.Range("J19").Value = 341
$ws.Range("K19").Value = -7.917888563049
$ws.Range("L19").Value = -2.484472049689
$ws.Range("M19").Value = 23.13725490196
$ws.Range("N19").Value = -20.906801007556

# --- Row 20 ---
$ws.Range("C20").Value = 3
$zeroRef.Copy($ws.Range("D20"))
$naRef.Copy($ws.Range("E20"))
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 58
$ws.Range("K20").Value = 34.883720930232
$ws.Range("L20").Value = 38.095238095238
$ws.Range("M20").Value = 222.222222222222
$ws.Range("N20").Value = -51.666666666666

# --- Row 21 ---
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 21.052631578947
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -18.681318681318
$ws.Range("I21").Value = 824
$ws.Range("J21").Value = 888
$ws.Range("K21").Value = -7.207207207207
$ws.Range("L21").Value = -4.959630911188
$ws.Range("M21").Value = 12.568306010929
$ws.Range("N21").Value = -68.669201520912

# --- Row 22 ---
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = 81.818181818181

# --- Row 23 ---
$numRef15.Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 74
$ws.Range("J23").Value = 77
$ws.Range("K23").Value = -3.896103896103
$ws.Range("L23").Value = 10.447761194029
$ws.Range("M23").Value = 94.736842105263

# --- Row 24 ---
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = -1.298701298701
$ws.Range("I24").Value = 1043
$ws.Range("J24").Value = 959
$ws.Range("K24").Value = 8.759124087591
$ws.Range("L24").Value = 4.195804195804
$ws.Range("M24").Value = 24.463007159904

# --- Row 25 ---
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 343
$ws.Range("J25").Value = 342
$ws.Range("K25").Value = 0.292397660818
$ws.Range("L25").Value = -0.867052023121
$ws.Range("M25").Value = -19.672131147541

# --- Row 26 ---
$zeroRef.Copy($ws.Range("D26"))
$naRef.Copy($ws.Range("E26"))
$ws.Range("L26").Value = -13.333333333333

# --- Row 27 ---
$numRef15.Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$numRef15.Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$numRef16.Copy($ws.Range("E27"))
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = -30.232558139534
$ws.Range("L27").Value = -25

# --- Row 28 ---
$zeroRef.Copy($ws.Range("D28"))
$naRef.Copy($ws.Range("E28"))
$ws.Range("G28").Value = 3
$ws.Range("M28").Value = -52.631578947368

# --- Row 29 ---
$zeroRef.Copy($ws.Range("D29"))
$naRef.Copy($ws.Range("E29"))
$ws.Range("G29").Value = 2
$ws.Range("M29").Value = -52.941176470588
